$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BOM")
$bom.Rows("6:12").Delete() | Out-Null

$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = '=BOM!$A$1:$AV$5'
$bom.PageSetup.PrintArea = 'A1:S5'

$bom.Range("C5").Select() | Out-Null

$nr = $wb.Worksheets.Item("NR")
$nr.Rows(1).RowHeight = 112.5
